$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8938  # 杭州·2024首届COMIC GALAXY次元盛典: 8928 -> 8938
$ws1.Range("F4").Value = 1966  # 杭州·浮游猫动漫嘉年华: 1964 -> 1966
$ws1.Range("F5").Value = 6592  # 杭州·理想乡动漫展-同人创作者大会: 6590 -> 6592
$ws1.Range("F7").Value = 2126  # 杭州·Eternal时光国乙only展（日+夜场）: 2125 -> 2126
$ws1.Range("F8").Value = 595  # 杭州·第五人格同人only: 593 -> 595
$ws1.Range("F9").Value = 74  # 杭州·鸳鸯谱婚配主题代号鸢同人only展: 72 -> 74
$ws1.Range("F13").Value = 5  # 桐庐·唯泽动漫游戏嘉年华: 4 -> 5
$ws1.Range("F16").Value = 8833  # 杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！: 8826 -> 8833
$ws1.Range("F17").Value = 166  # 杭州·第二届次元格子动漫展嘉宾内场——吴磊: 165 -> 166
$ws1.Range("F18").Value = 67  # 杭州·第二届次元格子动漫展嘉宾内场——赵乾景: 66 -> 67
$ws1.Range("F21").Value = 1835  # 杭州·第六届华盟次元动漫嘉年华: 1833 -> 1835
$ws1.Range("F25").Value = 77  # 杭州·弹丸论破only同人展: 74 -> 77
$ws1.Range("F27").Value = 200  # 杭州·第二届次元格子动漫展嘉宾内场——赵成晨: 199 -> 200
$ws1.Range("F29").Value = 8  # 杭州·逐月节·园游会·原神×绝区零×崩铁×崩坏同人only: 7 -> 8
$ws1.Range("F30").Value = 61  # 杭州·2024·华彩的摔跤宴 - 木吉KAZUYA降临: 58 -> 61
$ws1.Range("F31").Value = 429  # 杭州·文豪野犬同人only2.0: 309 -> 429
$ws1.Range("F33").Value = 15  # 杭州·火影同人ONLY: 12 -> 15
$ws1.Range("F34").Value = 419  # 杭州·第五人格同人only2.0: 301 -> 419
$ws1.Range("F35").Value = 2275  # 杭州·首届CCPC动漫嘉年华: 2267 -> 2275
$ws1.Range("F36").Value = 872  # 杭州·鸢飞鱼跃代号鸢only: 870 -> 872
$ws1.Range("F37").Value = 530  # 杭州·SK怀旧动漫展SK12: 524 -> 530
$ws1.Range("F41").Value = 277  # 杭州·亿万心动国乙✘代号鸢同人only(日夜场）: 269 -> 277
$ws1.Range("F44").Value = 1037  # 杭州·第六届AP动漫游戏嘉年华: 1035 -> 1037

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 212  # 杭州·2024CJMF·不止音乐节 卡琳娜专场: 211 -> 212
$ws2.Range("F14").Value = 9  # 杭州·东方乐典2024: 7 -> 9
$ws2.Range("F16").Value = 28  # 杭州·【逐漫旋律】跨越二次元经典动漫 ACG音乐会: 27 -> 28

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2340  # 杭州·盗墓笔记官方授权「四季同书」主题店: 2339 -> 2340
$ws3.Range("F3").Value = 719  # 杭州·剑网3×HAPPY ZOO 剑网3十五周年主题咖啡厅: 718 -> 719
$ws3.Range("F4").Value = 326  # 杭州·木灵动漫 二哈和他的白猫师尊主题餐厅: 325 -> 326

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2340  # 杭州·盗墓笔记官方授权「四季同书」主题店: 2339 -> 2340
$ws4.Range("F3").Value = 719  # 杭州·剑网3×HAPPY ZOO 剑网3十五周年主题咖啡厅: 718 -> 719
$ws4.Range("F4").Value = 212  # 杭州·2024CJMF·不止音乐节 卡琳娜专场: 211 -> 212
$ws4.Range("F6").Value = 8938  # 杭州·2024首届COMIC GALAXY次元盛典: 8928 -> 8938
$ws4.Range("F8").Value = 326  # 杭州·木灵动漫 二哈和他的白猫师尊主题餐厅: 325 -> 326
$ws4.Range("F9").Value = 1966  # 杭州·浮游猫动漫嘉年华: 1964 -> 1966
$ws4.Range("F10").Value = 2126  # 杭州·Eternal时光国乙only展（日+夜场）: 2125 -> 2126
$ws4.Range("F11").Value = 595  # 杭州·第五人格同人only: 593 -> 595
$ws4.Range("F12").Value = 74  # 杭州·鸳鸯谱婚配主题代号鸢同人only展: 72 -> 74
$ws4.Range("F20").Value = 8833  # 杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！: 8826 -> 8833
$ws4.Range("F21").Value = 166  # 杭州·第二届次元格子动漫展嘉宾内场——吴磊: 165 -> 166
$ws4.Range("F22").Value = 67  # 杭州·第二届次元格子动漫展嘉宾内场——赵乾景: 66 -> 67
$ws4.Range("F24").Value = 1835  # 杭州·第六届华盟次元动漫嘉年华: 1833 -> 1835
$ws4.Range("F28").Value = 77  # 杭州·弹丸论破only同人展: 75 -> 77
$ws4.Range("F29").Value = 200  # 杭州·第二届次元格子动漫展嘉宾内场——赵成晨: 199 -> 200
$ws4.Range("F32").Value = 429  # 杭州·文豪野犬同人only2.0: 310 -> 429
$ws4.Range("F34").Value = 419  # 杭州·第五人格同人only2.0: 301 -> 419
$ws4.Range("F35").Value = 2275  # 杭州·首届CCPC动漫嘉年华: 2267 -> 2275
$ws4.Range("F36").Value = 872  # 杭州·鸢飞鱼跃代号鸢only: 870 -> 872
$ws4.Range("F38").Value = 530  # 杭州·SK怀旧动漫展SK12: 524 -> 530
$ws4.Range("F39").Value = 277  # 杭州·亿万心动国乙✘代号鸢同人only(日夜场）: 269 -> 277
